$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.261.77"
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("D3").Value = "3.177.15"
$ws.Range("E3").Value = "  -8.17%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "564.62"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").Value = "169.01"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("D9").Value = "3.174.75"
$ws.Range("E9").Value = "  -8.13%  "
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  -6.33%  "
$ws.Range("D11").Value = "6.59"
$ws.Range("E11").Value = "  -5.41%  "
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  -5.35%  "
$ws.Range("D13").Value = "3.726.06"
$ws.Range("E13").Value = "  -8.15%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "27.35"
$ws.Range("E15").Value = "  -8.87%  "
$ws.Range("D16").Value = "64.260.63"
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("D18").Value = "3.175.62"
$ws.Range("E18").Value = "  -8.13%  "
$ws.Range("D19").Value = "5.73"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").Value = "353.27"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").Value = "7.19"
$ws.Range("E22").Value = "  -6.12%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "68.66"
$ws.Range("E24").Value = "  -6.44%  "
$ws.Range("D25").Value = "0.504"
$ws.Range("E25").Value = "  -6.02%  "
$ws.Range("E26").Value = "  -6.91%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "5.56"
$ws.Range("E31").Value = "  -6.65%  "
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").Value = "22.03"
$ws.Range("E33").Value = "  -7.01%  "
$ws.Range("D34").Value = "6.63"
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("D36").Value = "1.43"
$ws.Range("E36").Value = "  -8.50%  "
$ws.Range("D37").Value = "153.23"
$ws.Range("E37").Value = "  -5.09%  "
$ws.Range("D38").Value = "0.817"
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("D39").Value = "26.28"
$ws.Range("E39").Value = "  -6.40%  "
$ws.Range("E40").Value = "  -6.78%  "
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("D42").Value = "2.614.76"
$ws.Range("E42").Value = "  -7.52%  "
$ws.Range("E43").Value = "  -7.40%  "
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("D45").Value = "39.37"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "0.0642"
$ws.Range("E46").Value = "  -7.33%  "
$ws.Range("E47").Value = "  -6.27%  "
$ws.Range("D48").Value = "'319.50"
$ws.Range("E48").Value = "  -6.64%  "
$ws.Range("E49").Value = "  -6.98%  "
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  +0.01%  "
